# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (used only by the Notes Master)
#   ppt/theme/theme2.xml  -> "Integral" / "Red Violet" (used by the Slide Master,
#                             and therefore by every slide layout + slide)
#
# The target edit swaps the *content* of the two parts: the Slide Master's
# theme (theme2.xml) becomes the plain "Office Theme" colour scheme, while
# the Notes Master's theme (theme1.xml) becomes the former "Integral" / "Red
# Violet" colour scheme. Font scheme and format scheme (fills/lines/effects)
# are identical between the two themes already, so only the 12 colour-scheme
# slots need to change.
#
# Re-colour the presentation's (Slide Master's) theme to the "Office Theme"
# palette via the ThemeColorScheme COM surface (PowerPoint RGB long values
# are 0xBBGGRR, i.e. R + G*256 + B*65536).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# 1 dk1      -> 000000
$colors.Item(1).RGB = 0
# 2 lt1      -> FFFFFF
$colors.Item(2).RGB = 16777215
# 3 dk2      -> 44546A
$colors.Item(3).RGB = 6968388
# 4 lt2      -> E7E6E6
$colors.Item(4).RGB = 15132391
# 5 accent1  -> 5B9BD5
$colors.Item(5).RGB = 13998939
# 6 accent2  -> ED7D31
$colors.Item(6).RGB = 3243501
# 7 accent3  -> A5A5A5
$colors.Item(7).RGB = 10855845
# 8 accent4  -> FFC000
$colors.Item(8).RGB = 49407
# 9 accent5  -> 4472C4
$colors.Item(9).RGB = 12874308
# 10 accent6 -> 70AD47
$colors.Item(10).RGB = 4697456
# 11 hlink   -> 0563C1
$colors.Item(11).RGB = 12673797
# 12 folHlink -> 954F72
$colors.Item(12).RGB = 7491477
